# Update "Training Dashboard" sheet: decrement PERIOD TO EXPIRE (col H) by 1
# and bump LAST UPDATE (col I) from 03-Nov-2025 to 04-Nov-2025 for rows 3-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 25; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE
    $iCell.Value2 = "'04-Nov-2025"
}
